$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price observation was added; insert a row at 97, pushing the
# existing rows 97..124 down to 98..125.
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new record's data.
$ws.Range("A97").Value = 7
$ws.Range("B97").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C97").Value = "Ñuble"
$ws.Range("D97").Value = 45209
$ws.Range("E97").Value = 16
$ws.Range("F97").Value = 100112001
$ws.Range("G97").Value = "Berenjena"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 30
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = 10000
$ws.Range("N97").Value = "$/caja 60 unidades"
$ws.Range("O97").Value = "Región de Arica y Parinacota"
$ws.Range("P97").Value = 167
$ws.Range("Q97").Value = 60
$ws.Range("R97").Value = "Hortaliza"
